# Adiciona um novo bullet (nível 2) ao final do conteúdo do slide 6,
# informando a quantidade de notícias sumarizadas na amostra.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Insere um novo parágrafo após o texto existente, mantendo o mesmo
# nível de recuo (lvl 1 / IndentLevel 2) dos parágrafos irmãos.
$newPara = $tr.InsertAfter("`rSumarização de 287 notícias")
$newPara.IndentLevel = 2
